$d = $word.ActiveDocument

# Highlight color used across all metric highlights: hex 2C3E50 -> BGR decimal for Word's Font.Color
$metricColor = 5258796

function Bold-Terms($para, [string[]]$terms) {
    # Walk the supplied terms left-to-right inside this paragraph only,
    # bolding + coloring each occurrence in order so we never touch the
    # same run twice or bleed into neighboring paragraphs.
    $searchStart = $para.Range.Start
    $pEnd = $para.Range.End
    foreach ($term in $terms) {
        $rng = $d.Range($searchStart, $pEnd)
        $found = $rng.Find.Execute($term, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if ($found) {
            $rng.Font.Bold = 1
            $rng.Font.Color = $metricColor
            $searchStart = $rng.End
        }
    }
}

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text

    if ($t -like "*Discovered systematic race coding errors*") {
        Bold-Terms $p @("23%", "64%")
    }
    elseif ($t -like "*Utilized advanced sampling methods*") {
        Bold-Terms $p @([char]0x00B1 + "4.2%", [char]0x00B1 + "2.1%", "71%", "87%")
    }
    elseif ($t -like "*Trigonometric algorithm for boundary estimation*") {
        Bold-Terms $p @("73.5%", "`$4.7M")
    }
    elseif ($t -like "*Built real-time FEC analysis systems*") {
        Bold-Terms $p @("`$2")
    }
    elseif ($t -like "*Algorithmic innovation: Pioneered trigonometric*") {
        Bold-Terms $p @("73.5%")
    }
    elseif ($t -like "*`$4.7M savings enabled nonprofit access*") {
        Bold-Terms $p @("`$4.7M")
    }
    elseif ($t -like "*accuracy improvement in racial classification algorithms*") {
        Bold-Terms $p @("178%")
    }
}
